$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.732.23"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "3.303.83"
$ws.Range("E3").Value = "  +5.94%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.68"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.37"
$ws.Range("E6").Value = "  +4.84%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.302.13"
$ws.Range("E8").Value = "  +6.03%  "
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.48"
$ws.Range("E11").Value = "  +5.09%  "
$ws.Range("E12").Value = "  +4.19%  "
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.68"
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").Value = "3.847.76"
$ws.Range("E15").Value = "  +5.95%  "
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "3.303.08"
$ws.Range("E17").Value = "  +6.20%  "
$ws.Range("D18").Value = "63.805.58"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("E19").Value = "  +3.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.19"
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.14"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("E22").Value = "  +5.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.02"
$ws.Range("E23").Value = "  +5.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.53"
$ws.Range("E24").Value = "  +5.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.69"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.29"
$ws.Range("E28").Value = "  +7.31%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  +4.50%  "
$ws.Range("E31").Value = "  +4.77%  "
$ws.Range("E32").Value = "  +10.27%  "
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("E35").Value = "  +2.52%  "
$ws.Range("E36").Value = "  +4.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.11"
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("E38").Value = "  +8.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0401"
$ws.Range("E39").Value = "  +4.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "426.15"
$ws.Range("E40").Value = "  +2.78%  "
$ws.Range("D41").Value = "3.059.53"
$ws.Range("E41").Value = "  +5.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.38"
$ws.Range("E42").Value = "  +2.64%  "
$ws.Range("E43").Value = "  +3.55%  "
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("E46").Value = "  +5.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.30"
$ws.Range("E47").Value = "  +4.21%  "
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.96"
$ws.Range("E51").Value = "  +11.46%  "
